# Applies the "cryptos list" refresh described in the commit:
#   "Updated cryptos list on Wed Sep 11 11:39:16 UTC 2024 with GitHub Actions"
#
# Columns: A=index (unchanged), B=Coin, C=Link, D=Price, E=Volume(1h)
# D/E are stored as plain text in the workbook (e.g. "56.634.58", "  -0.06%  ").
# Many Price values look like plain numbers (e.g. "513.79"); Excel would silently
# coerce those into numeric cells, dropping trailing zeros / exact formatting.
# Set-TextValue forces such values to stay text, the same way a user gets text-kept
# numbers by prefixing a cell entry with a leading apostrophe in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    if ($Text -match '^[+-]?\d+(\.\d+)?$') {
        # Looks like a plain number -> force text with a leading apostrophe
        $Range.Value = "'" + $Text
    } else {
        $Range.Value = $Text
    }
}

# Row 2
Set-TextValue $ws.Range("D2") '56.634.58'
$ws.Range("E2").Value = '  -0.06%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.319.14'
$ws.Range("E3").Value = '  -0.24%  '

# Row 4
$ws.Range("E4").Value = '  -0.60%  '

# Row 5
Set-TextValue $ws.Range("D5") '513.79'
$ws.Range("E5").Value = '  -1.44%  '

# Row 6
Set-TextValue $ws.Range("D6") '131.80'
$ws.Range("E6").Value = '  -2.34%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.998'
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.534'
$ws.Range("E8").Value = '  -0.83%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.101'
$ws.Range("E9").Value = '  -3.16%  '

# Row 10
$ws.Range("E10").Value = '  -0.25%  '

# Row 11
Set-TextValue $ws.Range("D11") '5.24'
$ws.Range("E11").Value = '  -1.39%  '

# Row 12
$ws.Range("E12").Value = '  -1.92%  '

# Row 13
Set-TextValue $ws.Range("D13") '23.55'
$ws.Range("E13").Value = '  -1.83%  '

# Row 14
Set-TextValue $ws.Range("D14") '2.729.66'
$ws.Range("E14").Value = '  -1.00%  '

# Row 15
Set-TextValue $ws.Range("D15") '56.589.03'
$ws.Range("E15").Value = '  -0.39%  '

# Row 16
Set-TextValue $ws.Range("D16") '0.0000133'
$ws.Range("E16").Value = '  -1.24%  '

# Row 17
Set-TextValue $ws.Range("D17") '2.315.61'
$ws.Range("E17").Value = '  -0.70%  '

# Row 18
Set-TextValue $ws.Range("D18") '10.36'
$ws.Range("E18").Value = '  -1.41%  '

# Row 19
Set-TextValue $ws.Range("D19") '328.04'
$ws.Range("E19").Value = '  +1.55%  '

# Row 20
Set-TextValue $ws.Range("D20") '4.15'
$ws.Range("E20").Value = '  -1.94%  '

# Row 21
Set-TextValue $ws.Range("D21") '6.73'
$ws.Range("E21").Value = '  +2.01%  '

# Row 22
Set-TextValue $ws.Range("D22") '0.999'
$ws.Range("E22").Value = '  -0.19%  '

# Row 23
Set-TextValue $ws.Range("D23") '61.19'
$ws.Range("E23").Value = '  +0.45%  '

# Row 24
$ws.Range("E24").Value = '  -1.18%  '

# Row 25
Set-TextValue $ws.Range("D25") '8.59'
$ws.Range("E25").Value = '  +7.54%  '

# Row 26
Set-TextValue $ws.Range("D26") '0.999'
$ws.Range("E26").Value = '  +0.46%  '

# Row 27
$ws.Range("E27").Value = '  +1.09%  '

# Row 28
Set-TextValue $ws.Range("D28") '167.87'
$ws.Range("E28").Value = '  +0.21%  '

# Row 29
Set-TextValue $ws.Range("D29") '1.68'
$ws.Range("E29").Value = '  -2.48%  '

# Row 30
Set-TextValue $ws.Range("D30") '0.0₃0719'
$ws.Range("E30").Value = '  -3.37%  '

# Row 31
$ws.Range("E31").Value = '  -2.01%  '

# Row 32
Set-TextValue $ws.Range("D32") '18.28'
$ws.Range("E32").Value = '  -0.46%  '

# Row 34
Set-TextValue $ws.Range("D34") '0.996'
$ws.Range("E34").Value = '  +0.32%  '

# Row 35
$ws.Range("E35").Value = '  -1.08%  '

# Row 36
$ws.Range("E36").Value = '  -2.65%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.883'
$ws.Range("E37").Value = '  -5.06%  '

# Row 38
$ws.Range("E38").Value = '  +0.46%  '

# Row 39
Set-TextValue $ws.Range("D39") '38.60'
$ws.Range("E39").Value = '  +1.72%  '

# Row 40
Set-TextValue $ws.Range("D40") '149.04'
$ws.Range("E40").Value = '  +7.87%  '

# Row 41
$ws.Range("E41").Value = '  -1.69%  '

# Row 42
$ws.Range("E42").Value = '  -1.16%  '

# Row 43
Set-TextValue $ws.Range("D43") '275.93'
$ws.Range("E43").Value = '  -0.50%  '

# Row 44
$ws.Range("E44").Value = '  -3.95%  '

# Row 45
Set-TextValue $ws.Range("D45") '0.0927'
$ws.Range("E45").Value = '  -0.68%  '

# Row 46
Set-TextValue $ws.Range("D46") '0.0494'
$ws.Range("E46").Value = '  -2.41%  '

# Row 47
Set-TextValue $ws.Range("D47") '0.554'

# Row 48
Set-TextValue $ws.Range("D48") '18.22'
$ws.Range("E48").Value = '  +1.95%  '

# Row 49
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D49") '0.0214'
$ws.Range("E49").Value = '  -1.44%  '

# Row 50
$ws.Range("B50").Value = 'Polygon'
$ws.Range("C50").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range("D50") '0.378'
$ws.Range("E50").Value = '  -0.35%  '

# Row 51
Set-TextValue $ws.Range("D51") '17.08'
$ws.Range("E51").Value = '  +1.13%  '
